$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BAEPAbCiPC")
$ws.Activate()

# Update the boolean value for "hydrogen" row (row 22, column B) from 1 to 0
$ws.Range("B22").Value = 0

# Update the selected cell/range in the sheet view
$ws.Range("D26").Select()
